$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values
$ws.Range("D45").Value = 100
$ws.Range("D46").Value = 234

# Update selection / view to show cell D47 as active, scrolled so row 37 is at top
$ws.Range("D47").Select()
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
